$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.502.23"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.957.49"
$ws.Range("E3").Value = "  -0.82%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "583.21"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +7.49%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "157.10"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.29%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.746"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -2.01%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "53.47"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  -1.98%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "10.79"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "4.584.31"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "3.972.36"
$ws.Range("E15").Value = "  -0.43%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.28"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +7.63%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "13.96"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "72.280.56"
$ws.Range("E20").Value = "  +0.76%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "431.20"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.65"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +8.43%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "95.75"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.41"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.27%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "14.24"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "4.44"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +22.78%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "11.13"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.67"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +1.03%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "36.26"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.29%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.85"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +4.25%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "13.55"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.131"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "48.64"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "676.16"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "69.08"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.83%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.434"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "0.0₃0853"
$ws.Range("E39").Value = "  +0.05%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("E43").Value = "  +0.10%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.84"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +9.87%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0484"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E46").Value = "  -1.19%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.66"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.99%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "3.39"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +3.79%  "
$ws.Range("E50").Value = "  -0.29%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.14"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +6.44%  "
